$d = $word.ActiveDocument

# Replace every occurrence of $search with $replacement, without letting Word's
# AutoFormat/AutoCorrect "smart quotes" feature mangle straight apostrophes in
# the replacement text (Range.Text assignment bypasses that, unlike
# Find.Execute's built-in Replace parameter).
function Replace-All($search, $replacement) {
    $count = 0
    $rng = $d.Content
    while ($rng.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false)) {
        $rng.Text = $replacement
        $rng.Collapse(0)
        $rng.End = $d.Content.End
        $count++
    }
    return $count
}

# Long paragraph replacements first (search using the ORIGINAL text, before
# the job title / company name get renamed elsewhere in the document).
Replace-All "I am excited to apply for the Graduate Product Engineer position at Attio. The role aligns perfectly with my skills and aspirations, especially in revolutionizing business operations in the AI era, a field that strongly interests me. Attio's focus on contributing to key areas of the product and making product decisions resonates with my passion - my experience developing a full-stack food ordering platform has equipped me with skills in React, Node.js, and MySQL, and I am eager to contribute while growing with your team." "I am excited to apply for the Junior Front End Developer position at AllSaints. The role aligns perfectly with my skills and aspirations, especially in making customers feel cool and confident, a field that strongly interests me. AllSaints' focus on developing and optimizing website content pages resonates with my passion - having built a full-stack food ordering platform where I significantly improved customer experience through front-end optimizations, and I am eager to contribute while growing with your team."

Replace-All "I am a Full Stack Engineer who recently developed a full-stack food ordering platform for a café. This experience strengthened my experience in React, Node.js, and MySQL, and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Graduate Product Engineer position at Attio includes:" "I am a Full Stack Engineer who recently built a full-stack food ordering platform with real-time order processing. This experience strengthened my proficiency in HTML, CSS, JavaScript, and responsive design and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Junior Front End Developer position at AllSaints includes:"

Replace-All "My unique background as Full Stack Engineering | Cover Letter Generator has provided me with developed AI-powered applications using React and Node.js, which I believe can also contribute to driving the company’s success in achieving the company's goal." "My unique background as a Full Stack Engineer | Marketing Content Management Platform has provided me with experience in designing a multi-version content management system that enhances team collaboration, which I believe can also contribute to driving the company’s success in achieving the company's goal."

# Name change: "Amy Han Hsun Shih" -> "Amy Han Hsun Shi" (both occurrences)
Replace-All "Amy Han Hsun Shih" "Amy Han Hsun Shi"

# Job title in heading: "Graduate Product Engineer" -> "Junior Front End Developer"
Replace-All "Graduate Product Engineer" "Junior Front End Developer"

# Date: "8, Apr" -> "5, May"
Replace-All "8, Apr" "5, May"

# "To the hiring team at Attio" -> "To the hiring team at AllSaints"
Replace-All "To the hiring team at Attio" "To the hiring team at AllSaints"

# Bullet 1
Replace-All "Increased international customer engagement by 10%." "Improved customer experience by 10% through internationalization."

# Bullet 2
Replace-All "Integrated secure user authentication and detailed analytics." "Experience in optimizing front-end features for enhanced performance."

# Bullet 3
Replace-All "Gained valuable insight into user engagement and product improvement." "Real-time features significantly enhance user engagement."

# Email text change - remove leading 'a' (use Range.Text so the hyperlink
# element itself -- and its relationship -- stays intact).
$rngEmail = $d.Content
$rngEmail.Find.Execute("a23514788@gmail.com", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rngEmail.Text = "23514788@gmail.com"

# Phone number change: delete + insert (rather than a plain Range.Text
# assignment) so the preceding " " run is not merged into this one.
$rngPhone = $d.Content
$rngPhone.Find.Execute("07 366-318-764", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$rngPhone.Delete()
$rngPhone.InsertAfter("07 366318764")
